# F05 Froze Token Embeddings + Decoder 12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Epoch Accuracy (column B) values
$ws.Range("B3").Value = 0.28125
$ws.Range("B4").Value = 0.234375
$ws.Range("B6").Value = 0.203125
$ws.Range("B9").Value = 0.15625
$ws.Range("B10").Value = 0.15625
$ws.Range("B11").Value = 0.15625
$ws.Range("B16").Value = 0.15625
$ws.Range("B17").Value = 0.15625
$ws.Range("B22").Value = 0.15625
$ws.Range("B23").Value = 0.15625
$ws.Range("B24").Value = 0.15625
$ws.Range("B25").Value = 0.15625
$ws.Range("B26").Value = 0.15625
$ws.Range("B27").Value = 0.15625
$ws.Range("B28").Value = 0.15625
$ws.Range("B29").Value = 0.15625
$ws.Range("B30").Value = 0.15625
$ws.Range("B37").Value = 0.15625
$ws.Range("B38").Value = 0.15625
$ws.Range("B39").Value = 0.15625
$ws.Range("B40").Value = 0.125
$ws.Range("B41").Value = 0.125
$ws.Range("B42").Value = 0.125
$ws.Range("B43").Value = 0.125
$ws.Range("B44").Value = 0.125
$ws.Range("B45").Value = 0.125
$ws.Range("B46").Value = 0.125
$ws.Range("B47").Value = 0.125
$ws.Range("B48").Value = 0.125
$ws.Range("B49").Value = 0.125
$ws.Range("B50").Value = 0.125
$ws.Range("B51").Value = 0.125
$ws.Range("B52").Value = 0.125
$ws.Range("B53").Value = 0.125
$ws.Range("B54").Value = 0.125
$ws.Range("B55").Value = 0.125
$ws.Range("B56").Value = 0.125
$ws.Range("B57").Value = 0.125
$ws.Range("B58").Value = 0.125
$ws.Range("B59").Value = 0.125
$ws.Range("B60").Value = 0.125
$ws.Range("B61").Value = 0.125
$ws.Range("B62").Value = 0.125
$ws.Range("B63").Value = 0.125
$ws.Range("B64").Value = 0.125
$ws.Range("B65").Value = 0.125
$ws.Range("B66").Value = 0.125
$ws.Range("B67").Value = 0.125
$ws.Range("B68").Value = 0.125
$ws.Range("B69").Value = 0.125
$ws.Range("B70").Value = 0.125
$ws.Range("B71").Value = 0.125
$ws.Range("B72").Value = 0.125
$ws.Range("B73").Value = 0.125
$ws.Range("B74").Value = 0.125
$ws.Range("B75").Value = 0.125
$ws.Range("B76").Value = 0.125
$ws.Range("B104").Value = 0.09375
$ws.Range("B105").Value = 0.171875
$ws.Range("B106").Value = 0.125
$ws.Range("B107").Value = 0.0625
$ws.Range("B109").Value = 0.09375
$ws.Range("B110").Value = 0.171875
$ws.Range("B113").Value = 0.15625
$ws.Range("B116").Value = 0.0625

# Update column A DisplayOutputs memory address repr for rows 102-118
for ($r = 102; $r -le 118; $r++) {
    $ws.Cells.Item($r, 1).Value = "<__main__.DisplayOutputs object at 0x7f3948cae8e0>"
}
